$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the localization status flips from
# "Ready for handoff" to "In Translation" everywhere it appears (the
# Overview sheet's per-language status columns, plus each language
# sheet's own Status column), and the Status column narrows to fit the
# shorter text.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: columns E (zh-cn) and F (de-de) hold the per-language status.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-language sheets: column C is "Status".
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the Status column(s) to fit the new, shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
